# Actualización 25 de Marzo
# Appends the new daily-report row (row 23) to the Chile-by-region case
# counts sheet, and moves the sheet's viewport/selection the way the
# author's session left it after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row ----------------------------------------------
# Copy the formatting of the last existing row (22) down into the new
# row (23) first, so the date cell (column A) keeps its DD/MM/YY format,
# then fill in the actual values for 24-Mar-2020 (day 22 of the series).
$ws.Range("A22:S22").Copy() | Out-Null
$ws.Range("A23:S23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$rowValues = @(43914, 22, 2, 4, 19, 1, 12, 32, 682, 11, 30, 111, 95, 74, 11, 44, 1, 13, 1142)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(23, $col).Value = $rowValues[$col - 1]
}

# --- Move viewport / selection ------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 8
$ws.Range("H26:Y30").Select() | Out-Null
